$d = $word.ActiveDocument

# The first two paragraphs currently hold the italic title run(s) + a line
# break + "Chapter 1 - Why ==============================", and a bold
# "By Dorothy Day" byline. They need to become a single Title-styled
# paragraph ("From Union Square to Rome Chapter 1 - Why"), split word by
# word into its own run (as produced by a pandoc title-block conversion),
# followed by a plain paragraph "% Dorothy Day".

$p1 = $d.Paragraphs.Item(1)
$p2 = $d.Paragraphs.Item(2)
$targetRange = $d.Range($p1.Range.Start, $p2.Range.End)

$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

$titleWords = @("From", "Union", "Square", "to", "Rome", "Chapter", "1", "-", "Why")

$runsXml = ""
for ($i = 0; $i -lt $titleWords.Length; $i++) {
    if ($i -gt 0) {
        $runsXml += "<w:r><w:t xml:space=`"preserve`"> </w:t></w:r>"
    }
    $runsXml += "<w:r><w:t xml:space=`"preserve`">$($titleWords[$i])</w:t></w:r>"
}

$titleParagraphXml = "<w:p $wNs><w:pPr><w:pStyle w:val=`"Title`"/></w:pPr>$runsXml</w:p>"
$bylineParagraphXml = "<w:p $wNs><w:r><w:t xml:space=`"preserve`">% Dorothy Day</w:t></w:r></w:p>"

$null = $targetRange.InsertXML($titleParagraphXml + $bylineParagraphXml)
